# "No need for any data in control panel"
# Models will now run with no data included in control panel.
#
# Clears the input values from the control_panel sheet (keeping the
# cell styles/formatting intact) and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")

# Clear single-cell inputs (values only, formatting/styles untouched)
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()

# Row 10 "reference" values (D through J) get cleared too
$ws.Range("D10:J10").ClearContents()

# Remaining single-value inputs further down the sheet
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("B20").ClearContents()

# Update the recorded selection to match the new active cell
$ws.Range("B11").Select() | Out-Null
